$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old content first (old sheet used A1:B6)
$ws.Range("A1:C23").Clear()

# --- Teacher block ---
$ws.Range("A1").Value = "Teacher"
$ws.Range("A2").Value = "a"
$ws.Range("A3").Value = "b"

# --- Student block ---
$ws.Range("A5").Value = "Student"
$ws.Range("A6").Value = "aa"
$ws.Range("A7").Value = "bb"

# --- Grades for student "aa" ---
$ws.Range("A9").Value = "aa"
$ws.Range("B9").Value = "math"
$ws.Range("C9").Value = 2
$ws.Range("B10").Value = "hist"
$ws.Range("C10").Value = 99
$ws.Range("B11").Value = "Average:"
$ws.Range("C11").Value = 50.5

# --- Grades for student "bb" ---
$ws.Range("A13").Value = "bb"
$ws.Range("B13").Value = "math"
$ws.Range("C13").Value = 93
$ws.Range("B14").Value = "hist"
$ws.Range("C14").Value = 84
$ws.Range("B15").Value = "Average:"
$ws.Range("C15").Value = 88.5

# --- Students summary block ---
$ws.Range("A18").Value = "Students Average:"
$ws.Range("B18").Value = 69.5
$ws.Range("A19").Value = "Median:"
$ws.Range("B19").Value = 69.5
$ws.Range("A20").Value = "Excellent Students:"
$ws.Range("A21").Value = "bb"
$ws.Range("B21").Value = 88.5

# --- Extra row ---
$ws.Range("A23").Value = "bb"
$ws.Range("B23").Value = "b"
$ws.Range("C23").Value = 88.5
